$d = $word.ActiveDocument
$t = $d.Tables(1)

# Row 1 (1-based Word index): "99.97" -> "0M"
$t.Cell(1,1).Range.Text = "0M"

# Row 2: "0.12" -> "0M"
$t.Cell(2,1).Range.Text = "0M"

# Row 3: "427" -> "0M"
$t.Cell(3,1).Range.Text = "0M"

# Row 4: "201" -> "804"
$t.Cell(4,1).Range.Text = "804"

# Row 6: "0.00017" -> "0.00061"
$t.Cell(6,1).Range.Text = "0.00061"

# Row 7: "0.00007" -> "0.00015"
$t.Cell(7,1).Range.Text = "0.00015"

# Row 8: "0.00002" -> "0.00003"
$t.Cell(8,1).Range.Text = "0.00003"

# Row 9: "0.00006" -> "0.00022"
$t.Cell(9,1).Range.Text = "0.00022"

# Row 10: "0.00007" -> "0.00025"
$t.Cell(10,1).Range.Text = "0.00025"

# Row 11: "0.00009" -> "0.00031"
$t.Cell(11,1).Range.Text = "0.00031"

# Row 12: "0.01496" -> "0.12192"
$t.Cell(12,1).Range.Text = "0.12192"

# Row 44: collapsed tab-separated run -> "99.97"
$t.Cell(44,1).Range.Text = "99.97"

# Row 45: collapsed tab-separated run -> "0.12"
$t.Cell(45,1).Range.Text = "0.12"

# Row 46: collapsed tab-separated run -> "427"
$t.Cell(46,1).Range.Text = "427"
